# "update complete, cancel genectics case"
# - Sheet1!A2 (the "id" column of the genetics case row) is updated to a new
#   generated case id: CA-R723QKZS -> CA-5EDNCJRQ
# - Sheet2!B6 (the "Status" column for the "assign follow up to approved 3"
#   test case row) is switched from pass -> fail, i.e. the case is cancelled.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws1.Range("A2").Value = "CA-5EDNCJRQ"
$ws2.Range("B6").Value = "fail"
